{"js": "// This script replaces the date heading and all the multiplication\n// problems in the practice-sheet table with the day's new values.\n// Every \"old\" search string in this document is unique, so a direct\n// search-and-replace for each pair is safe and order independent.\nconst replacements = [\n  [\"2025-11-10 Monday\", \"2025-11-11 Tuesday\"],\n  [\"302\u00d75=\", \"384\u00d79=\"],\n  [\"823\u00d72=\", \"392\u00d74=\"],\n  [\"663\u00d76=\", \"322\u00d76=\"],\n  [\"187\u00d74=\", \"834\u00d76=\"],\n  [\"377\u00d72=\", \"934\u00d74=\"],\n  [\"235\u00d76=\", \"183\u00d76=\"],\n  [\"331\u00d79=\", \"752\u00d73=\"],\n  [\"515\u00d73=\", \"259\u00d72=\"],\n  [\"553\u00d79=\", \"436\u00d73=\"],\n  [\"657\u00d79=\", \"781\u00d78=\"],\n  [\"829\u00d76=\", \"724\u00d74=\"],\n  [\"520\u00d74=\", \"269\u00d77=\"],\n  [\"504\u00d76=\", \"745\u00d78=\"],\n  [\"516\u00d79=\", \"358\u00d79=\"],\n  [\"780\u00d78=\", \"311\u00d78=\"],\n  [\"276\u00d76=\", \"515\u00d74=\"],\n  [\"635\u00d74=\", \"257\u00d78=\"],\n  [\"424\u00d73=\", \"342\u00d78=\"],\n  [\"337\u00d73=\", \"547\u00d77=\"],\n  [\"758\u00d77=\", \"139\u00d74=\"],\n  [\"401\u00d74=\", \"109\u00d74=\"],\n  [\"704\u00d75=\", \"486\u00d75=\"],\n  [\"480\u00d73=\", \"271\u00d73=\"],\n  [\"116\u00d73=\", \"568\u00d76=\"],\n  [\"603\u00d77=\", \"110\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Updates the date heading and every multiplication problem in the\n# practice-sheet table to the day's new values. Every \"old\" string below\n# is unique within the document, so a straight Find/Replace per pair,\n# applied to the whole document Range, is safe regardless of order.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"2025-11-10 Monday\", \"2025-11-11 Tuesday\"),\n    @(\"302\u00d75=\", \"384\u00d79=\"),\n    @(\"823\u00d72=\", \"392\u00d74=\"),\n    @(\"663\u00d76=\", \"322\u00d76=\"),\n    @(\"187\u00d74=\", \"834\u00d76=\"),\n    @(\"377\u00d72=\", \"934\u00d74=\"),\n    @(\"235\u00d76=\", \"183\u00d76=\"),\n    @(\"331\u00d79=\", \"752\u00d73=\"),\n    @(\"515\u00d73=\", \"259\u00d72=\"),\n    @(\"553\u00d79=\", \"436\u00d73=\"),\n    @(\"657\u00d79=\", \"781\u00d78=\"),\n    @(\"829\u00d76=\", \"724\u00d74=\"),\n    @(\"520\u00d74=\", \"269\u00d77=\"),\n    @(\"504\u00d76=\", \"745\u00d78=\"),\n    @(\"516\u00d79=\", \"358\u00d79=\"),\n    @(\"780\u00d78=\", \"311\u00d78=\"),\n    @(\"276\u00d76=\", \"515\u00d74=\"),\n    @(\"635\u00d74=\", \"257\u00d78=\"),\n    @(\"424\u00d73=\", \"342\u00d78=\"),\n    @(\"337\u00d73=\", \"547\u00d77=\"),\n    @(\"758\u00d77=\", \"139\u00d74=\"),\n    @(\"401\u00d74=\", \"109\u00d74=\"),\n    @(\"704\u00d75=\", \"486\u00d75=\"),\n    @(\"480\u00d73=\", \"271\u00d73=\"),\n    @(\"116\u00d73=\", \"568\u00d76=\"),\n    @(\"603\u00d77=\", \"110\u00d78=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
